# "inetbanking updated adding new customer"
# The login test-data sheet holds username/password pairs. The oldest
# customer record (row 4: mngr1111 / jttttt) is retired and row 3's
# record (mngr164225 / jahetAp) is replaced by a newly added customer
# (mngr523220 / gynUnYd), leaving a fresh blank row behind it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new customer's credentials in place of the old row 3 entry.
$ws.Range("A3").Value = "mngr523220"
$ws.Range("B3").Value = "gynUnYd"

# Drop the stale row 4 customer entirely, keeping the row as blank.
$ws.Range("A4:B5").Value = $null

# Touch the vacated rows' formatting so they stay part of the sheet's
# used range as blank rows (matches the still-present row 4/5 stubs).
$ws.Rows.Item(4).Font.Size = 11
$ws.Rows.Item(5).Font.Size = 11

# Selection ends up on the newly-edited record.
$ws.Range("A3").Select()
